$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The LR-pair table used to include every Sending x Target cluster combination,
# including "self" pairs where Sending cluster == Target cluster. The refreshed
# TPM-based run drops those self pairs, so only 6 of the original 9 data rows
# remain, and every remaining metric column is recomputed against the new TPM values.
$data = New-Object 'object[,]' 6,20
$data[0,0] = "ECs"
$data[0,1] = "Col1a1"
$data[0,2] = "Itga11"
$data[0,3] = "FAPs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 1.120168333333333
$data[0,7] = 3.360505
$data[0,8] = 0.001768092629909379
$data[0,9] = 0.001768092629909379
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 33.94639966666666
$data[0,13] = 101.839199
$data[0,14] = 0.9746097333921855
$data[0,15] = 0.9746097333921855
$data[0,16] = 38.02568193727721
$data[0,17] = 342.231137435495
$data[0,18] = 0.001723200286648668
$data[0,19] = 0.001723200286648668
$data[1,0] = "ECs"
$data[1,1] = "Col1a1"
$data[1,2] = "Itga11"
$data[1,3] = "MuSCs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 1.120168333333333
$data[1,7] = 3.360505
$data[1,8] = 0.001768092629909379
$data[1,9] = 0.001768092629909379
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.8843623333333334
$data[1,13] = 2.653087
$data[1,14] = 0.02539026660781448
$data[1,15] = 0.02539026660781448
$data[1,16] = 0.9906346809927777
$data[1,17] = 8.915712128935
$data[1,18] = 0.00004489234326071098
$data[1,19] = 0.00004489234326071098
$data[2,0] = "FAPs"
$data[2,1] = "Col1a1"
$data[2,2] = "Itga11"
$data[2,3] = "FAPs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 604.0312093333333
$data[2,7] = 1812.093628
$data[2,8] = 0.9534130698726969
$data[2,9] = 0.9534130698726969
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 33.94639966666666
$data[2,13] = 101.839199
$data[2,14] = 0.9746097333921855
$data[2,15] = 0.9746097333921855
$data[2,16] = 20504.68484316933
$data[2,17] = 184542.163588524
$data[2,18] = 0.9292056578412543
$data[2,19] = 0.9292056578412543
$data[3,0] = "FAPs"
$data[3,1] = "Col1a1"
$data[3,2] = "Itga11"
$data[3,3] = "MuSCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 604.0312093333333
$data[3,7] = 1812.093628
$data[3,8] = 0.9534130698726969
$data[3,9] = 0.9534130698726969
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.8843623333333334
$data[3,13] = 2.653087
$data[3,14] = 0.02539026660781448
$data[3,15] = 0.02539026660781448
$data[3,16] = 534.1824496921818
$data[3,17] = 4807.642047229637
$data[3,18] = 0.02420741203144263
$data[3,19] = 0.02420741203144263
$data[4,0] = "MuSCs"
$data[4,1] = "Col1a1"
$data[4,2] = "Itga11"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 28.39480333333333
$data[4,7] = 85.18441
$data[4,8] = 0.04481883749739363
$data[4,9] = 0.04481883749739363
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 33.94639966666666
$data[4,13] = 101.839199
$data[4,14] = 0.9746097333921855
$data[4,15] = 0.9746097333921855
$data[4,16] = 963.9013424097321
$data[4,17] = 8675.11208168759
$data[4,18] = 0.04368087526428249
$data[4,19] = 0.04368087526428249
$data[5,0] = "MuSCs"
$data[5,1] = "Col1a1"
$data[5,2] = "Itga11"
$data[5,3] = "MuSCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 28.39480333333333
$data[5,7] = 85.18441
$data[5,8] = 0.04481883749739363
$data[5,9] = 0.04481883749739363
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.8843623333333334
$data[5,13] = 2.653087
$data[5,14] = 0.02539026660781448
$data[5,15] = 0.02539026660781448
$data[5,16] = 25.11129453040778
$data[5,17] = 226.00165077367
$data[5,18] = 0.001137962233111137
$data[5,19] = 0.001137962233111137

$ws.Range("A2:T7").Value = $data

# Remove the old trailing rows (previously rows 8-10, the MuSCs-sending self-pair
# block) so the used range shrinks from A1:T10 down to A1:T7.
$ws.Range("A8:T10").Delete()
